$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.844.94"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.637.08"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.511"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.866.39"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.652.59"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "26.835.95"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.66%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +3.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  +4.45%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "1.258.91"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.808"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").Value = "1.781.01"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E51").Value = "  -0.33%  "
